$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 753 (header)
$ws.Range("A69:F69").Copy()
$ws.Range("A753:F753").PasteSpecial(-4122)
$ws.Cells.Item(753, "C").Value = 'TUESDAY'

# Row 754 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A754:F754").PasteSpecial(-4122)
$ws.Cells.Item(754, "A").Value = 'Demo'
$ws.Cells.Item(754, "B").Value = 42745
$ws.Cells.Item(754, "C").Value = '1550'
$ws.Cells.Item(754, "D").Value = 'R'
$ws.Cells.Item(754, "E").Value = 'S203'
$ws.Cells.Item(754, "F").Value = 'Meet Professor Lawrence Lam.'

# Row 755 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A755:F755").PasteSpecial(-4122)
$ws.Cells.Item(755, "A").Value = 'Pickup PC'
$ws.Cells.Item(755, "B").Value = 42745
$ws.Cells.Item(755, "C").Value = '1630'
$ws.Cells.Item(755, "D").Value = 'VC'
$ws.Cells.Item(755, "E").Value = '113'
$ws.Cells.Item(755, "F").Value = 'Flat screen DLP and wireless keyboard. To Vanier 132 storeroom.'

# Row 756 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A756:F756").PasteSpecial(-4122)
$ws.Cells.Item(756, "A").Value = 'AV Shutdown'
$ws.Cells.Item(756, "B").Value = 42745
$ws.Cells.Item(756, "C").Value = '1730'
$ws.Cells.Item(756, "D").Value = 'MC'
$ws.Cells.Item(756, "E").Value = '101A'
$ws.Cells.Item(756, "F").Value = 'Pick up wireless keyboard and TV remote control. To FDRS 164.'

# Row 757 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A757:F757").PasteSpecial(-4122)
$ws.Cells.Item(757, "A").Value = 'AV Shutdown'
$ws.Cells.Item(757, "B").Value = 42745
$ws.Cells.Item(757, "C").Value = '1900'
$ws.Cells.Item(757, "D").Value = 'R'
$ws.Cells.Item(757, "E").Value = 'N203'

# Row 758 (alt)
$ws.Range("A74:F74").Copy()
$ws.Range("A758:F758").PasteSpecial(-4122)
$ws.Cells.Item(758, "A").Value = 'Other'
$ws.Cells.Item(758, "B").Value = 42745
$ws.Cells.Item(758, "C").Value = '2030'
$ws.Cells.Item(758, "D").Value = 'MC'
$ws.Cells.Item(758, "E").Value = '157A'
$ws.Cells.Item(758, "F").Value = 'Door code 11012* '

# Row 759 (header)
$ws.Range("A69:F69").Copy()
$ws.Range("A759:F759").PasteSpecial(-4122)
$ws.Cells.Item(759, "C").Value = 'WEDNESDAY'

# Row 760 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A760:F760").PasteSpecial(-4122)
$ws.Cells.Item(760, "A").Value = 'Pickup PC'
$ws.Cells.Item(760, "B").Value = 42746
$ws.Cells.Item(760, "C").Value = '1600'
$ws.Cells.Item(760, "D").Value = 'CSQ'
$ws.Cells.Item(760, "E").Value = 'East Bear Pit'
$ws.Cells.Item(760, "F").Value = '3 flat screen TVs, wireless keyboards (turn off), 2 TV remotes and extension cords. Two of the PC''s have wireless network receivers attached, the other is hard wired to the wall jack.  Return all to Vari 1155 and connect all 3 PC''s to wired internet router there.'
$ws.Rows.Item(760).RowHeight = 60

# Row 761 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A761:F761").PasteSpecial(-4122)
$ws.Cells.Item(761, "A").Value = 'Pickup Skype Kit'
$ws.Cells.Item(761, "B").Value = 42746
$ws.Cells.Item(761, "C").Value = '1700'
$ws.Cells.Item(761, "D").Value = 'R'
$ws.Cells.Item(761, "E").Value = 'N940'
$ws.Cells.Item(761, "F").Value = 'Web cam and tripod to Ross S120.'

# Row 762 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A762:F762").PasteSpecial(-4122)
$ws.Cells.Item(762, "A").Value = 'AV Shutdown'
$ws.Cells.Item(762, "B").Value = 42746
$ws.Cells.Item(762, "C").Value = '1700'
$ws.Cells.Item(762, "D").Value = 'R'
$ws.Cells.Item(762, "E").Value = 'N940'
$ws.Cells.Item(762, "F").Value = '3 wired mics, neck mic, stands, projector remote and presentation remote to back booth.'
$ws.Rows.Item(762).RowHeight = 30

# Row 763 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A763:F763").PasteSpecial(-4122)
$ws.Cells.Item(763, "A").Value = 'Demo'
$ws.Cells.Item(763, "B").Value = 42746
$ws.Cells.Item(763, "C").Value = '1850'
$ws.Cells.Item(763, "D").Value = 'R'
$ws.Cells.Item(763, "E").Value = 'N203'
$ws.Cells.Item(763, "F").Value = 'Meet instructor who shall remain nameless for the moment.'

# Row 764 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A764:F764").PasteSpecial(-4122)
$ws.Cells.Item(764, "A").Value = 'Demo'
$ws.Cells.Item(764, "B").Value = 42746
$ws.Cells.Item(764, "C").Value = '1850'
$ws.Cells.Item(764, "D").Value = 'VH'
$ws.Cells.Item(764, "E").Value = '3009'
$ws.Cells.Item(764, "F").Value = 'Meet instructor who shall remain nameless for the moment.'

# Row 765 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A765:F765").PasteSpecial(-4122)
$ws.Cells.Item(765, "A").Value = 'AV Shutdown'
$ws.Cells.Item(765, "B").Value = 42746
$ws.Cells.Item(765, "C").Value = '1900'
$ws.Cells.Item(765, "D").Value = 'R'
$ws.Cells.Item(765, "E").Value = 'N102'
$ws.Cells.Item(765, "F").Value = 'Nat Taylor Cinema. Lock cinema all doors after shutdown.'

# Row 766 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A766:F766").PasteSpecial(-4122)
$ws.Cells.Item(766, "A").Value = 'AV Shutdown'
$ws.Cells.Item(766, "B").Value = 42746
$ws.Cells.Item(766, "C").Value = '2030'
$ws.Cells.Item(766, "D").Value = 'R'
$ws.Cells.Item(766, "E").Value = 'N203'

# Row 767 (alt)
$ws.Range("A74:F74").Copy()
$ws.Range("A767:F767").PasteSpecial(-4122)
$ws.Cells.Item(767, "A").Value = 'Other'
$ws.Cells.Item(767, "B").Value = 42746
$ws.Cells.Item(767, "C").Value = '2030'
$ws.Cells.Item(767, "D").Value = 'MC'
$ws.Cells.Item(767, "E").Value = '157A'
$ws.Cells.Item(767, "F").Value = 'Door code 11012* '

# Row 768 (header)
$ws.Range("A69:F69").Copy()
$ws.Range("A768:F768").PasteSpecial(-4122)
$ws.Cells.Item(768, "C").Value = 'THURSDAY'

# Row 769 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A769:F769").PasteSpecial(-4122)
$ws.Cells.Item(769, "A").Value = 'Pickup PC'
$ws.Cells.Item(769, "B").Value = 42747
$ws.Cells.Item(769, "C").Value = '1630'
$ws.Cells.Item(769, "D").Value = 'VC'
$ws.Cells.Item(769, "E").Value = '113'
$ws.Cells.Item(769, "F").Value = 'Flat screen DLP and wireless keyboard. To Vanier 132 storeroom.'

# Row 770 (alt)
$ws.Range("A74:F74").Copy()
$ws.Range("A770:F770").PasteSpecial(-4122)
$ws.Cells.Item(770, "A").Value = 'Other'
$ws.Cells.Item(770, "B").Value = 42747
$ws.Cells.Item(770, "C").Value = '1730'
$ws.Cells.Item(770, "D").Value = 'MC'
$ws.Cells.Item(770, "E").Value = '157A'
$ws.Cells.Item(770, "F").Value = 'Door code 11012* '

# Row 771 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A771:F771").PasteSpecial(-4122)
$ws.Cells.Item(771, "A").Value = 'AV Shutdown'
$ws.Cells.Item(771, "B").Value = 42747
$ws.Cells.Item(771, "C").Value = '1730'
$ws.Cells.Item(771, "D").Value = 'R'
$ws.Cells.Item(771, "E").Value = 'N102'
$ws.Cells.Item(771, "F").Value = 'Nat Taylor Cinema. Lock cinema all doors after shutdown.'

# Row 772 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A772:F772").PasteSpecial(-4122)
$ws.Cells.Item(772, "A").Value = 'AV Shutdown'
$ws.Cells.Item(772, "B").Value = 42747
$ws.Cells.Item(772, "C").Value = '1800'
$ws.Cells.Item(772, "D").Value = 'R'
$ws.Cells.Item(772, "E").Value = 'N940'
$ws.Cells.Item(772, "F").Value = 'Return lectern mic, 2 desk mics, stands, projector remote and presentation remote to back booth.'
$ws.Rows.Item(772).RowHeight = 30

# Row 773 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A773:F773").PasteSpecial(-4122)
$ws.Cells.Item(773, "A").Value = 'AV Shutdown'
$ws.Cells.Item(773, "B").Value = 42747
$ws.Cells.Item(773, "C").Value = '1900'
$ws.Cells.Item(773, "D").Value = 'R'
$ws.Cells.Item(773, "E").Value = 'N203'

# Row 774 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A774:F774").PasteSpecial(-4122)
$ws.Cells.Item(774, "A").Value = 'AV Shutdown'
$ws.Cells.Item(774, "B").Value = 42747
$ws.Cells.Item(774, "C").Value = '1900'
$ws.Cells.Item(774, "D").Value = 'R'
$ws.Cells.Item(774, "E").Value = 'S203'

# Row 775 (header)
$ws.Range("A69:F69").Copy()
$ws.Range("A775:F775").PasteSpecial(-4122)
$ws.Cells.Item(775, "C").Value = 'FRIDAY'

# Row 776 (alt)
$ws.Range("A74:F74").Copy()
$ws.Range("A776:F776").PasteSpecial(-4122)
$ws.Cells.Item(776, "A").Value = 'Other'
$ws.Cells.Item(776, "B").Value = 42748
$ws.Cells.Item(776, "C").Value = '1730'
$ws.Cells.Item(776, "D").Value = 'MC'
$ws.Cells.Item(776, "E").Value = '157A'
$ws.Cells.Item(776, "F").Value = 'Door code 11012* '

# Row 777 (header)
$ws.Range("A69:F69").Copy()
$ws.Range("A777:F777").PasteSpecial(-4122)
$ws.Cells.Item(777, "C").Value = 'MONDAY'

# Row 778 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A778:F778").PasteSpecial(-4122)
$ws.Cells.Item(778, "A").Value = 'Demo'
$ws.Cells.Item(778, "B").Value = 42751
$ws.Cells.Item(778, "C").Value = '1620'
$ws.Cells.Item(778, "D").Value = 'MC'
$ws.Cells.Item(778, "E").Value = '216'
$ws.Cells.Item(778, "F").Value = 'Meet instructor Asma Sidddiqi'

# Row 779 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A779:F779").PasteSpecial(-4122)
$ws.Cells.Item(779, "A").Value = 'AV Shutdown'
$ws.Cells.Item(779, "B").Value = 42751
$ws.Cells.Item(779, "C").Value = '1630'
$ws.Cells.Item(779, "D").Value = 'MC'
$ws.Cells.Item(779, "E").Value = '101A'
$ws.Cells.Item(779, "F").Value = 'Pick up wireless keyboard and TV remote control. To FDRS 164.'

# Row 780 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A780:F780").PasteSpecial(-4122)
$ws.Cells.Item(780, "A").Value = 'AV Shutdown'
$ws.Cells.Item(780, "B").Value = 42751
$ws.Cells.Item(780, "C").Value = '1630'
$ws.Cells.Item(780, "D").Value = 'R'
$ws.Cells.Item(780, "E").Value = 'N940'
$ws.Cells.Item(780, "F").Value = 'No mics were used. Pc and projector only.'

# Row 781 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A781:F781").PasteSpecial(-4122)
$ws.Cells.Item(781, "A").Value = 'Pickup PC'
$ws.Cells.Item(781, "B").Value = 42751
$ws.Cells.Item(781, "C").Value = '1630'
$ws.Cells.Item(781, "D").Value = 'VC'
$ws.Cells.Item(781, "E").Value = '221'
$ws.Cells.Item(781, "F").Value = 'Flat screen DLP and wireless keyboard. To Vanier 132 storeroom.'

# Row 782 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A782:F782").PasteSpecial(-4122)
$ws.Cells.Item(782, "A").Value = 'AV Shutdown'
$ws.Cells.Item(782, "B").Value = 42751
$ws.Cells.Item(782, "C").Value = '1730'
$ws.Cells.Item(782, "D").Value = 'R'
$ws.Cells.Item(782, "E").Value = 'N203'

# Row 783 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A783:F783").PasteSpecial(-4122)
$ws.Cells.Item(783, "A").Value = 'AV Shutdown'
$ws.Cells.Item(783, "B").Value = 42751
$ws.Cells.Item(783, "C").Value = '1830'
$ws.Cells.Item(783, "D").Value = 'R'
$ws.Cells.Item(783, "E").Value = 'S203'

# Row 784 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A784:F784").PasteSpecial(-4122)
$ws.Cells.Item(784, "A").Value = 'AV Shutdown'
$ws.Cells.Item(784, "B").Value = 42751
$ws.Cells.Item(784, "C").Value = '1830'
$ws.Cells.Item(784, "D").Value = 'R'
$ws.Cells.Item(784, "E").Value = 'N102'
$ws.Cells.Item(784, "F").Value = 'Nat Taylor Cinema. Lock all cinema doors after shutdown.'

# Row 785 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A785:F785").PasteSpecial(-4122)
$ws.Cells.Item(785, "A").Value = 'Demo'
$ws.Cells.Item(785, "B").Value = 42751
$ws.Cells.Item(785, "C").Value = '1850'
$ws.Cells.Item(785, "D").Value = 'VH'
$ws.Cells.Item(785, "E").Value = 'D'
$ws.Cells.Item(785, "F").Value = 'Meet instructor Carole Bigwood.  Monitor cutting out?'

# Row 786 (header)
$ws.Range("A69:F69").Copy()
$ws.Range("A786:F786").PasteSpecial(-4122)
$ws.Cells.Item(786, "C").Value = 'TUESDAY'

# Row 787 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A787:F787").PasteSpecial(-4122)
$ws.Cells.Item(787, "A").Value = 'AV Shutdown'
$ws.Cells.Item(787, "B").Value = 42752
$ws.Cells.Item(787, "C").Value = '1730'
$ws.Cells.Item(787, "D").Value = 'MC'
$ws.Cells.Item(787, "E").Value = '101A'
$ws.Cells.Item(787, "F").Value = 'Pick up wireless keyboard and TV remote control. To FDRS 164.'

# Row 788 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A788:F788").PasteSpecial(-4122)
$ws.Cells.Item(788, "A").Value = 'AV Shutdown'
$ws.Cells.Item(788, "B").Value = 42752
$ws.Cells.Item(788, "C").Value = '1900'
$ws.Cells.Item(788, "D").Value = 'R'
$ws.Cells.Item(788, "E").Value = 'N203'

# Row 789 (alt)
$ws.Range("A74:F74").Copy()
$ws.Range("A789:F789").PasteSpecial(-4122)
$ws.Cells.Item(789, "A").Value = 'Other'
$ws.Cells.Item(789, "B").Value = 42752
$ws.Cells.Item(789, "C").Value = '2030'
$ws.Cells.Item(789, "D").Value = 'MC'
$ws.Cells.Item(789, "E").Value = '157A'
$ws.Cells.Item(789, "F").Value = 'Door code 11012* '

# Row 790 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A790:F790").PasteSpecial(-4122)
$ws.Cells.Item(790, "A").Value = 'Pickup Mic'
$ws.Cells.Item(790, "B").Value = 42752
$ws.Cells.Item(790, "C").Value = '1530'
$ws.Cells.Item(790, "D").Value = 'MC'
$ws.Cells.Item(790, "E").Value = '014 JCR'
$ws.Cells.Item(790, "F").Value = 'Pick up 3 desk mics, all cables and clips and NECK MIC ALSO. Return to Van 040 basement storeroom.'
$ws.Rows.Item(790).RowHeight = 30

# Row 791 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A791:F791").PasteSpecial(-4122)
$ws.Cells.Item(791, "A").Value = 'Pickup Large PA'
$ws.Cells.Item(791, "B").Value = 42752
$ws.Cells.Item(791, "C").Value = '1530'
$ws.Cells.Item(791, "D").Value = 'MC'
$ws.Cells.Item(791, "E").Value = '014 JCR'
$ws.Cells.Item(791, "F").Value = 'Pick up amplifier, 2 speaker cables and 2 large speakers. Return to Van 040 basement storeroom.'
$ws.Rows.Item(791).RowHeight = 30

# Row 792 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A792:F792").PasteSpecial(-4122)
$ws.Cells.Item(792, "A").Value = 'Other'
$ws.Cells.Item(792, "B").Value = 42752
$ws.Cells.Item(792, "C").Value = '1530'
$ws.Cells.Item(792, "D").Value = 'MC'
$ws.Cells.Item(792, "E").Value = '014 JCR'
$ws.Cells.Item(792, "F").Value = 'Pick up all matts and ac cords and return to Van 040 basement. Key for room in Fdrs 164 storeroom keyrack. '
$ws.Rows.Item(792).RowHeight = 30

# Row 793 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A793:F793").PasteSpecial(-4122)
$ws.Cells.Item(793, "A").Value = 'Demo'
$ws.Cells.Item(793, "B").Value = 42752
$ws.Cells.Item(793, "C").Value = '1630'
$ws.Cells.Item(793, "D").Value = 'VH'
$ws.Cells.Item(793, "E").Value = 'A'
$ws.Cells.Item(793, "F").Value = 'Built in PC not working - roll in PC cart in room going thru document camera input. Press "Doc cam" to "Projector to get image on screen. Demo to prof.'
$ws.Rows.Item(793).RowHeight = 45

# Row 794 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A794:F794").PasteSpecial(-4122)
$ws.Cells.Item(794, "A").Value = 'Demo'
$ws.Cells.Item(794, "B").Value = 42752
$ws.Cells.Item(794, "C").Value = '1900'
$ws.Cells.Item(794, "D").Value = 'VH'
$ws.Cells.Item(794, "E").Value = 'A'
$ws.Cells.Item(794, "F").Value = 'Built in PC not working - roll in PC cart in room going thru document camera input. Press "Doc cam" to "Projector to get image on screen. Demo to prof.'
$ws.Rows.Item(794).RowHeight = 45

# Row 795 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A795:F795").PasteSpecial(-4122)
$ws.Cells.Item(795, "A").Value = 'Pickup PC'
$ws.Cells.Item(795, "B").Value = 42752
$ws.Cells.Item(795, "C").Value = '2200'
$ws.Cells.Item(795, "D").Value = 'VH'
$ws.Cells.Item(795, "E").Value = 'A'
$ws.Cells.Item(795, "F").Value = 'Pick up roll in PC and return to Vari 1019 storeroom.'

# Row 796 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A796:F796").PasteSpecial(-4122)
$ws.Cells.Item(796, "A").Value = 'Other'
$ws.Cells.Item(796, "B").Value = 42752
$ws.Cells.Item(796, "C").Value = '1745'
$ws.Cells.Item(796, "D").Value = 'R'
$ws.Cells.Item(796, "E").Value = 'N102'
$ws.Cells.Item(796, "F").Value = 'Open up Nat Taylor cinema.'

# Row 797 (standard)
$ws.Range("A75:F75").Copy()
$ws.Range("A797:F797").PasteSpecial(-4122)
$ws.Cells.Item(797, "A").Value = 'AV Shutdown'
$ws.Cells.Item(797, "B").Value = 42752
$ws.Cells.Item(797, "C").Value = '2200'
$ws.Cells.Item(797, "D").Value = 'R'
$ws.Cells.Item(797, "E").Value = 'N102'
$ws.Cells.Item(797, "F").Value = 'Nat Taylor Cinema. Lock cinema all doors after shutdown.'

$excel.CutCopyMode = $false
$ws.Range("F800").Select()
Write-Host "done"